# Summary of uploaded dataset now shown on preview.html: rearrange two
# rows of sample data, drop the unused E:G "placeholder" columns, turn on
# an AutoFilter over the header row, set the page to portrait, and leave
# the selection where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Oranges" (row 4) and "Berries" (row 7) records ---
$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$c4 = $ws.Range("C4").Value()
$d4 = $ws.Range("D4").Value()

$a7 = $ws.Range("A7").Value()
$b7 = $ws.Range("B7").Value()
$c7 = $ws.Range("C7").Value()
$d7 = $ws.Range("D7").Value()

$ws.Range("A4").Value = $a7
$ws.Range("B4").Value = $b7
$ws.Range("C4").Value = $c7
$ws.Range("D4").Value = $d7

$ws.Range("A7").Value = $a4
$ws.Range("B7").Value = $b4
$ws.Range("C7").Value = $c4
$ws.Range("D7").Value = $d4

# --- Drop the stray empty E:G columns so the used range shrinks to A:D ---
$ws.Range("E1:G11").Clear()

# --- Turn on AutoFilter over the header row (adds the hidden
#     _xlnm._FilterDatabase defined name scoped to this sheet) ---
$flt = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$1")
$flt.Visible = $false

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Leave the selection on N6, matching where the user clicked last ---
[void]$ws.Range("N6").Select()

Write-Output "done"
